$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date bump in A1
$ws.Range("A1").Value = 45309

# Step 1 / Step 2 pricing updates (column D)
$ws.Range("D26").Value = 1246
$ws.Range("D27").Value = 1246
$ws.Range("D28").Value = 1467
$ws.Range("D29").Value = 2536
$ws.Range("D30").Value = 3128
$ws.Range("D31").Value = 4150
$ws.Range("D32").Value = 5903
$ws.Range("D33").Value = 8740

$ws.Range("D35").Value = 1840
$ws.Range("D36").Value = 2082
$ws.Range("D37").Value = 2657
$ws.Range("D38").Value = 4414
$ws.Range("D39").Value = 5725
$ws.Range("D40").Value = 7299
$ws.Range("D41").Value = 10263
$ws.Range("D42").Value = 15834
